# Generate Report for Handoff
#
# The handoff/report generation re-ran after the original "low"-priority
# rows were re-classified as "ht" and their handoff timestamps were
# refreshed:
#   - zh-cn!E4:E7   "low" -> "ht"                       (Priority)
#   - zh-cn!H4:H7   "2016-08-31 00:33:00" -> "2016-08-31 00:33:28"   (Latest Handoff Datetime)
#   - de-de!E4:E7   "low" -> "ht"                       (Priority)
#   - de-de!H4:H7 / Overview!G4:G7
#       "2016-08-31 00:33:12" -> "2016-08-31 00:33:32"  (shared Latest Handoff Datetime / Latest HO Xliff Generate Date)

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-31 00:33:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-31 00:33:32"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-31 00:33:32"
